# bom.xlsx: mark two pending orders as received ("Shipping" -> "On Hand")
# and record the pickup date for the two most recent shipments on the
# "Shipping Numbers" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Materials")
$ws2 = $wb.Worksheets.Item("Shipping Numbers")

# --- Materials sheet: flip Status column from "Shipping" to "On Hand" ---
$ws1.Range("G9").Value  = "On Hand"
$ws1.Range("G10").Value = "On Hand"
$ws1.Range("G31").Value = "On Hand"
$ws1.Range("G32").Value = "On Hand"

$ws1.Range("C36").Select()

# --- Shipping Numbers sheet: log "Date Picked Up" (column F) for rows 9-10 ---
$ws2.Range("F9").Value  = 44152
$ws2.Range("F10").Value = 44152

# Match the existing date formatting used elsewhere in column D/E/F
$ws2.Range("D9").Copy()
$ws2.Range("F9:F10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A12").Select()
